# SAO_BORJA.xlsx maintenance edit:
#  - rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - remove the (now unused) "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null
$excel.DisplayAlerts = $true
